$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "ToDo" text in column C for rows 13-21 (shift content down by
#    one logical slot, clearing the old C13 and editing a couple of entries
#    from "theory" -> "methods" wording as well as adding new rows at the end)
# ---------------------------------------------------------------------------

# Row 13 no longer has a ToDo entry
$ws.Range("C13").Value = ""

$ws.Range("C14").Value = "LSTM, Write Evaluation, put figures in, check what to change in figures, create all chapter outlines and write rought summary for all of them"
$ws.Range("C15").Value = "LSTM, Clean up evaluation/generate missing plots, identify what methods are needed to explain evaluation, write methods"
$ws.Range("C16").Value = "LSTM, Write most of methods, LSTM, metrics/loss function"
$ws.Range("C17").Value = "LSTM, Write creation of datasets, odes, normalization"
$ws.Range("C18").Value = "LSTM, Write training theory, hyperparameter overview, hyperparameter search manually vs autotuner"
$ws.Range("C19").Value = "Kann ich das glauben Seminar, go over presentation again"
$ws.Range("C20").Value = 'LSTM, Write "other work" chapter, introduction'
$ws.Range("C21").Value = "LSTM, Clean up writing, check for things that need more work, make new plan for rest"

# ---------------------------------------------------------------------------
# 2. Append new schedule rows 28-36 (dates 2022-12-01 .. 2022-12-09) with the
#    same Date/Weekday formula pattern used by the existing rows, copying the
#    number format / styling from the last existing row (27).
# ---------------------------------------------------------------------------

$ws.Range("A27:B27").Copy() | Out-Null
$ws.Range("A28:B36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$newDates = @(44896,44897,44898,44899,44900,44901,44902,44903,44904)
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = 28 + $i
    $ws.Cells.Item($row, 1).Value = $newDates[$i]
    $ws.Cells.Item($row, 2).FormulaArray = '=_xlfn.SWITCH(WEEKDAY(A' + $row + ',1),1,"Sun",2,"Mon",3,"Tue",4,"Wed",5,"Thu",6,"Fri",7,"Sat")'
}

$ws.Range("C28").Value = "Abgabe fast fertig"
$ws.Range("C34").Value = "Abgabe final"

# ---------------------------------------------------------------------------
# 3. Cosmetic sheet-level updates reflected in the diff
# ---------------------------------------------------------------------------

# Column C got a bit wider to fit the new text (closest achievable value via
# the COM ColumnWidth property, which is quantized to 1/6-character units)
$ws.Columns.Item(3).ColumnWidth = 126.42

# Selection moved to C24 (one row further down than before)
$ws.Range("C24").Select() | Out-Null
